# Revert "Update in sample.pptx"
#
# 1. Restore the cached "datetimeFigureOut" field text from 2019/12/11
#    back to 2019/5/8 on the slide master and on every slide layout.
# 2. Remove the "C8 c8 c8 " run that had been typed into the title
#    placeholder on slide 1.

$p = $ppt.ActivePresentation

$oldDate = "2019/12/11"
$newDate = "2019/5/8"

function Update-DateShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -eq -1) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master
$master = $p.SlideMaster
Update-DateShape $master.Shapes

# Every slide layout attached to the master
for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DateShape $layout.Shapes
}

# Slide 1: drop the "C8 c8 c8 " text that was added to the title placeholder
$s1 = $p.Slides.Item(1)
for ($i = 1; $i -le $s1.Shapes.Count; $i++) {
    $shp = $s1.Shapes.Item($i)
    if ($shp.HasTextFrame -eq -1) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "C8 c8 c8 ") {
            $tr.Text = ""
        }
    }
}
